$wb = $excel.ActiveWorkbook

# New device-name text shared by all four summary sheets.
$newDeviceText = "NGC-1854/T772 OR TC-65368 "

# Per-sheet "Input Value" (B8) updates. GreenColorPercentage (sheetId 1) keeps its
# original value of 10, so it is omitted from this update map.
$b8Updates = @{
    "VDWorstCaseYellowPercentage"  = 20
    "VtgDropYellowColorPercentage" = 35
    "RedColorPercentage"           = 37
}

foreach ($ws in $wb.Worksheets) {
    # Update the "Device Name" value in B4 and drop its border formatting
    # (the cell reverts to the workbook's default/unstyled look).
    $ws.Range("B4").Borders.LineStyle = -4142   # xlLineStyleNone
    $ws.Range("B4").Value = $newDeviceText
    $ws.Rows.Item(4).RowHeight = 28.8

    # Move the active selection to B4, matching the saved view state.
    $ws.Range("B4").Select()

    if ($b8Updates.ContainsKey($ws.Name)) {
        $ws.Range("B8").Value = $b8Updates[$ws.Name]
    }
}
